$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(16.0711673480879,8.61511508945847,13.56159526462787,14.13141669706078,0,36.10223932728334,16.16056733970392,25.08492602862891,8.618319071558375,0,12.14483238649772,16.75017320725858,18.85690058318922,25.57501843936498),
    @(15.69554943976792,8.408031093794085,13.5774981841238,14.16852733367672,0,36.1458093438882,16.20508514765344,25.1826004031837,8.619227528606347,0,12.14922044457661,16.67288204849944,18.90594102950857,25.64019010945213),
    @(15.46227949579755,8.277174150467411,13.58953482735986,14.19274197816834,0,36.18270240426539,16.23494989178043,25.24702233676351,8.619895881063302,0,12.15335348796407,16.62731078251428,18.93778529807597,25.68538390738985),
    @(15.3666896081625,8.22296596002367,13.59501148646372,14.2029695594283,0,36.200281388143,16.24775623644609,25.27439324910999,8.620196143789517,0,12.15540047605049,16.6092272638776,18.951198880816,25.70510084451449),
    @(15.35078893978659,8.213912902161333,13.59595541193497,14.20468960006714,0,36.20335388419868,16.24992114220424,25.27900570857112,8.620247690534745,0,12.15576231177286,16.60625429877741,18.95345260484356,25.70845328830526),
    @(15.46099230506832,8.276446587563404,13.58960637272508,14.19287845268765,0,36.18292918483655,16.23512002678503,25.24738694241788,8.619899817393238,0,12.15337962443982,16.62706491276293,18.93796442855416,25.68564455561105),
    @(15.94228008341629,8.544510294960029,13.56660719585564,14.14391636339271,0,36.11515474553445,16.17539177627642,25.117680420891,8.61860942336447,0,12.14604736589889,16.72313973563197,18.87345056533739,25.59641376529624),
    @(16.85960512591909,9.038779619310173,13.53951894076091,14.05920547048267,0,36.06289292143094,16.07834966017473,24.89864890584472,8.616951756811378,0,12.14304100011405,16.92591466555052,18.760648595303,25.46260171937324),
    @(17.51048371539804,9.380337961984038,13.53057143878374,14.00381575625788,0,36.07382395965622,16.01930381855903,24.75928458459926,8.616260183216713,0,12.14770343296701,17.0828342683387,18.68607108290903,25.38949298055759),
    @(17.80025308898944,9.530584530576224,13.52887107297606,13.98009520722745,0,36.08951393350959,15.99510363960971,24.70057108443902,8.616058640964225,0,12.15130181565875,17.15576163848213,18.65393326485425,25.36172518599571),
    @(17.90897296383899,9.586707774367836,13.52856698138714,13.97132449112875,0,36.0969944002354,15.98632225460425,24.67901199504718,8.615998476394537,0,12.15287547078652,17.18358323980259,18.64201969741148,25.35200083297552),
    @(17.88560461457591,9.574655470899106,13.52861737762753,13.97320401366749,0,36.09531495004065,15.98819645979053,24.68362512059351,8.616010717040645,0,12.15252719419549,17.17758248096931,18.64457410912157,25.35405996132345),
    @(17.80921832047658,9.53521747808232,13.52883925128524,13.97936939481802,0,36.09009852952338,15.99437351900263,24.6987838833012,8.616053368003399,0,12.15142706041671,17.15804651310994,18.65294799564304,25.36090930407414),
    @(17.76229518244089,9.510959085693246,13.52901937535318,13.98317342667641,0,36.08710365702617,15.99820698960123,24.70815692150669,8.61608159377983,0,12.15078063470142,17.14610645029935,18.65811060181847,25.36520773221612),
    @(17.49141065591919,9.370412774817726,13.53073016910087,14.00539561136362,0,36.07301409670791,16.02093890081761,24.76321597941609,8.616275620731056,0,12.14749787978503,17.07809796310199,18.68820727541525,25.39141825431101),
    @(17.32353916093694,9.282852952366193,13.53238601866828,14.01940595662133,0,36.06711455435344,16.03556563844559,24.79819324074462,8.616423533918272,0,12.14586143186378,17.03676079499164,18.70712797418892,25.40890460439933),
    @(17.22639471808638,9.232009841094241,13.53356155435019,14.02760334685565,0,36.06473028912972,16.0442289202168,24.81875206524418,8.616519257358094,0,12.14505934271794,17.01313106723576,18.71817898674435,25.41947894751161),
    @(17.19340528617978,9.214713713423182,13.53399792175469,14.0304027352621,0,36.06409635390036,16.04720515371628,24.82578861625635,8.616553499543876,0,12.14481171067884,17.00515607291207,18.72194960370809,25.42314793171476),
    @(17.34147110375059,9.292223885637311,13.53218666371908,14.0179001486426,0,36.06763815405552,16.0339826835717,24.79442422109963,8.616406687054452,0,12.14602124283568,17.04114618170108,18.70509641789061,25.40698967043332),
    @(17.83168298946433,9.546822578307228,13.52876486817681,13.97755273090808,0,36.09158897485079,15.99254877935623,24.69431307604043,8.616040402753429,0,12.1517444809886,17.16377925305752,18.65048143156458,25.35887601733696),
    @(18.14613936551016,9.708704740413882,13.52850857452359,13.95241728682732,0,36.11621070300588,15.96770004707678,24.63281622431097,8.615895131569912,0,12.15671424465311,17.24511757762024,18.61628110708762,25.33204025117863),
    @(17.97888157258616,9.622728621276947,13.52846455501846,13.96571983130753,0,36.10225010654759,15.98075811897687,24.66527821633643,8.615964087702672,0,12.15394980086903,17.20160234515274,18.63439803618809,25.34594088772786),
    @(17.33336603593182,9.287988853344677,13.53227609564444,14.01858048002939,0,36.06739829631757,16.03469754573573,24.79612679343159,8.616414270228304,0,12.14594856007086,17.03916312485886,18.70601434476347,25.40785378864442),
    @(16.61502370944143,8.908692826075031,13.54492083496705,14.08091642091217,0,36.06837501807894,16.10245112444327,24.95411876424632,8.617307356894969,0,12.15130181565875,16.86960292942553,18.78970304886875,25.49438234559911),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
